$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Step 1: mint a brand-new numbered-list definition (numId) -------------
# Word allocates a new w:num / w:abstractNum pair in numbering.xml the first
# time a paragraph is given "1. 2. 3." numbering. We trigger that allocation
# on a throw-away scratch paragraph appended at the very end of the story,
# then delete the scratch paragraph again; the freshly minted list
# definition stays behind in numbering.xml for us to reference explicitly
# (numId=4) on the real list paragraphs inserted afterwards.
$scratchRange = $d.Content
$scratchRange.Collapse(0)
$scratchRange.InsertXML("<w:p $wNs><w:r><w:t>scratch</w:t></w:r></w:p>")
$scratchPara = $d.Paragraphs.Last
$scratchPara.Range.ListFormat.ApplyNumberDefault()
$scratchPara.Range.Delete()

# --- Step 2: append the new paragraphs at the end of the document ----------
$body = ""

# Blank spacer paragraph.
$body += "<w:p $wNs><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr></w:p>"

# "Question:" paragraph.
$body += "<w:p $wNs><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" `
       + "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Question:</w:t></w:r></w:p>"

# First numbered question.
$body += "<w:p $wNs><w:pPr><w:pStyle w:val='ListParagraph'/>" `
       + "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='4'/></w:numPr>" `
       + "<w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" `
       + "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr>" `
       + "<w:t>TCP ISN với TCP ID là các thông số được sử dụng bởi protocol để xác định các packet, segment thì làm sao có thể xác thông tin của OS được?</w:t></w:r></w:p>"

# Second numbered question.
$body += "<w:p $wNs><w:pPr><w:pStyle w:val='ListParagraph'/>" `
       + "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='4'/></w:numPr>" `
       + "<w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" `
       + "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr>" `
       + "<w:t>Cơ chế RPC của Nmap cụ thể là thế nào?</w:t></w:r></w:p>"

$insertRange = $d.Content
$insertRange.Collapse(0)
$insertRange.InsertXML($body)

Write-Host "Paragraphs now:" $d.Paragraphs.Count
